$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178748488426208
$ws.Range("B1").Value = 2.41413688659668
$ws.Range("D1").Value = 2.335357427597046
$ws.Range("E1").Value = 1.195580840110779
